$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 139.3
$ws.Range("I6").Value = 125.57895
$ws.Range("K6").Value = 376.73685
$ws.Range("M6").Value = -264.73685

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 36
$ws.Range("I8").Value = 36
$ws.Range("K8").Value = 108
$ws.Range("M8").Value = 31

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 135.88461
$ws.Range("I33").Value = 79
$ws.Range("K33").Value = 79
$ws.Range("M33").Value = 150

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1574.6086
$ws.Range("I100").Value = 1393.625
$ws.Range("J100").Value = 1988.2858
$ws.Range("K100").Value = 1393.625
$ws.Range("L100").Value = 1988.2858
$ws.Range("M100").Value = -852.625
$ws.Range("N100").Value = -3070.2858

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2715.125
$ws.Range("J113").Value = 3130.6365
$ws.Range("L113").Value = 3130.6365
$ws.Range("N113").Value = -9638.636500000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4979.1377
$ws.Range("I141").Value = 3667.1428
$ws.Range("J141").Value = 6203.6665
$ws.Range("K141").Value = 11001.4284
$ws.Range("L141").Value = 18610.9995
$ws.Range("M141").Value = -5821.428400000001
$ws.Range("N141").Value = -28970.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2579.6155
$ws.Range("I2").Value = 2924.6667
$ws.Range("J2").Value = 1803.25
$ws.Range("K2").Value = 2924.6667
$ws.Range("L2").Value = 1803.25
$ws.Range("M2").Value = -2811.6667
$ws.Range("N2").Value = -2029.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 9266.333000000001
$ws.Range("I12").Value = 2999
$ws.Range("J12").Value = 12400
$ws.Range("K12").Value = 2999
$ws.Range("L12").Value = 12400
$ws.Range("M12").Value = -2826
$ws.Range("N12").Value = -12746

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19912.225
$ws.Range("I32").Value = 22602.082
$ws.Range("K32").Value = 22602.082
$ws.Range("M32").Value = -22315.082

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2579.6155
$ws.Range("I116").Value = 2924.6667
$ws.Range("J116").Value = 1803.25
$ws.Range("K116").Value = 2924.6667
$ws.Range("L116").Value = 1803.25
$ws.Range("M116").Value = -630.6667000000002
$ws.Range("N116").Value = -6391.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2579.6155
$ws.Range("I3").Value = 2924.6667
$ws.Range("J3").Value = 1803.25
$ws.Range("K3").Value = 2924.6667
$ws.Range("L3").Value = 1803.25
$ws.Range("M3").Value = -2810.6667
$ws.Range("N3").Value = -2031.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 1000
$ws.Range("I8").Value = 1000
$ws.Range("K8").Value = 1000
$ws.Range("M8").Value = -860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 65780
$ws.Range("J50").Value = 65780
$ws.Range("L50").Value = 65780
$ws.Range("N50").Value = -66928

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 30000
$ws.Range("J57").Value = 30000
$ws.Range("L57").Value = 30000
$ws.Range("N57").Value = -31440

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1441.8182
$ws.Range("I94").Value = 1427.3846
$ws.Range("J94").Value = 1462.6666
$ws.Range("K94").Value = 1427.3846
$ws.Range("L94").Value = 1462.6666
$ws.Range("M94").Value = -976.3846000000001
$ws.Range("N94").Value = -2364.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 55960
$ws.Range("J124").Value = 55960
$ws.Range("L124").Value = 55960
$ws.Range("N124").Value = -65780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 30000
$ws.Range("J136").Value = 30000
$ws.Range("L136").Value = 30000
$ws.Range("N136").Value = -40200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 9800
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 9800
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 9800
$ws.Range("N3").Value = -10026
$ws.Range("M3").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 40000
$ws.Range("J137").Value = 40000
$ws.Range("L137").Value = 40000
$ws.Range("N137").Value = -50200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 3725.889
$ws.Range("J106").Value = 3725.889
$ws.Range("L106").Value = 11177.667
$ws.Range("N106").Value = -13069.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2611820
$ws.Range("I139").Value = 4404392
$ws.Range("J139").Value = 4443
$ws.Range("K139").Value = 13213176
$ws.Range("L139").Value = 13329
$ws.Range("M139").Value = -13208036
$ws.Range("N139").Value = -23609

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2138.48
$ws.Range("I140").Value = 1994
$ws.Range("J140").Value = 2716.4
$ws.Range("K140").Value = 5982
$ws.Range("L140").Value = 8149.200000000001
$ws.Range("M140").Value = -802
$ws.Range("N140").Value = -18509.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3650
$ws.Range("J113").Value = 2600
$ws.Range("L113").Value = 2600
$ws.Range("N113").Value = -6940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2933.7273
$ws.Range("I132").Value = 3357.182
$ws.Range("J132").Value = 2510.2727
$ws.Range("K132").Value = 10071.546
$ws.Range("L132").Value = 7530.8181
$ws.Range("M132").Value = -7541.545999999998
$ws.Range("N132").Value = -12590.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 44322.5
$ws.Range("J137").Value = 44322.5
$ws.Range("L137").Value = 44322.5
$ws.Range("N137").Value = -54522.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 500499.5
$ws.Range("J2").Value = 999999
$ws.Range("L2").Value = 999999
$ws.Range("N2").Value = -1000223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 388.06668
$ws.Range("I55").Value = 446
$ws.Range("J55").Value = 359.1
$ws.Range("K55").Value = 446
$ws.Range("L55").Value = 359.1
$ws.Range("M55").Value = -273
$ws.Range("N55").Value = -705.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 38792.668
$ws.Range("J87").Value = 38792.668
$ws.Range("L87").Value = 38792.668
$ws.Range("N87").Value = -41038.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 37766.75
$ws.Range("I88").Value = 30500
$ws.Range("K88").Value = 30500
$ws.Range("M88").Value = -30072

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H90").Value = 38792.668
$ws.Range("J90").Value = 38792.668
$ws.Range("L90").Value = 116378.004
$ws.Range("N90").Value = -127610.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H91").Value = 37766.75
$ws.Range("I91").Value = 30500
$ws.Range("K91").Value = 30500
$ws.Range("M91").Value = -29018

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3781.61
$ws.Range("I136").Value = 2253.3635
$ws.Range("K136").Value = 6760.0905
$ws.Range("M136").Value = -4210.0905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 1669999.6
$ws.Range("I3").Value = 1669999.6
$ws.Range("K3").Value = 1669999.6
$ws.Range("M3").Value = -1669885.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 19833.334
$ws.Range("I4").Value = 53666.668
$ws.Range("J4").Value = 2916.6667
$ws.Range("K4").Value = 53666.668
$ws.Range("L4").Value = 2916.6667
$ws.Range("M4").Value = -53553.668
$ws.Range("N4").Value = -3142.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 3438.4
$ws.Range("I6").Value = 525
$ws.Range("K6").Value = 525
$ws.Range("M6").Value = -410

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 22600
$ws.Range("J58").Value = 28900
$ws.Range("L58").Value = 28900
$ws.Range("N58").Value = -29516

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 33203.75
$ws.Range("J70").Value = 33203.75
$ws.Range("L70").Value = 33203.75
$ws.Range("N70").Value = -33833.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 33203.75
$ws.Range("J73").Value = 33203.75
$ws.Range("L73").Value = 33203.75
$ws.Range("N73").Value = -35387.75
